# Rebuild the worksheet data to reflect the new "combined table" output:
#  - Columns C and D are swapped: C now holds total_requested_2020,
#    D now holds total_returned_2016 (matching the new header order).
#  - Rows are re-sorted/re-grouped: states with both 2016 and 2020 data
#    are kept together in alphabetical order; two new states (HI, IN)
#    that only have 2016 data are inserted in alphabetical position;
#    four states that only have 2020 data (CT, OK, PA, SD) are moved
#    to the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("state", "total_requested_2016", "total_requested_2020", "total_returned_2016", "total_returned_2020", "diff_requested", "pctchg_requested", "diff_returned", "pctchg_returned"),
    @("AK", 45042, 129533, 28687, 49655, 84491, 187.58, 20968, 73.09),
    @("AR", 266795, 318602, 266780, 277162, 51807, 19.42, 10382, 3.89),
    @("AZ", 2340508, 3284338, 995460, 1310818, 943830, 40.33, 315358, 31.68),
    @("CA", 11178435, 5080324, 2480093, 5080324, -6098111, -54.55, 2600231, 104.84),
    @("CO", 3167702, 1271757, 813678, 1271757, -1895945, -59.85, 458079, 56.3),
    @("DE", 21245, 173228, 15703, 103597, 151983, 715.38, 87894, 559.73),
    @("FL", 4654451, 6696516, 3217937, 4173722, 2042065, 43.87, 955785, 29.7),
    @("GA", 1298213, 2647203, 1198262, 1870473, 1348990, 103.91, 672211, 56.1),
    @("HI", 46, $null, 0, $null, $null, $null, $null, $null),
    @("IA", 540549, 820240, 405546, 609019, 279691, 51.74, 203473, 50.17),
    @("ID", 144138, 428347, 112291, 233679, 284209, 197.18, 121388, 108.1),
    @("IL", 960015, 2343725, 697072, 874521, 1383710, 144.13, 177449, 25.46),
    @("IN", 414548, $null, 375812, $null, $null, $null, $null, $null),
    @("KS", 302252, 566120, 204584, 228913, 263868, 87.3, 24329, 11.89),
    @("LA", 296152, 502791, 296124, 502791, 206639, 69.77, 206667, 69.79),
    @("MD", 456737, 1634824, 357964, 799532, 1178087, 257.94, 441568, 123.36),
    @("ME", 189054, 402352, 136150, 271863, 213298, 112.82, 135713, 99.68),
    @("MI", 1066154, 3054940, 659190, 1925701, 1988786, 186.54, 1266511, 192.13),
    @("MN", 276138, 1121512, 276105, 1121512, 845374, 306.14, 845407, 306.19),
    @("MT", 325466, 632967, 178201, 255748, 307501, 94.48, 77547, 43.52),
    @("NC", 1548337, 2413084, 1443495, 2401839, 864747, 55.85, 958344, 66.39),
    @("ND", 80380, 204251, 50147, 125191, 123871, 154.11, 75044, 149.65),
    @("NE", 198095, 501734, 125727, 326478, 303639, 153.28, 200751, 159.67),
    @("NJ", 383216, 5898132, 226234, 1994520, 5514916, 1439.11, 1768286, 781.62),
    @("NM", 298475, 594845, 270349, 450933, 296370, 99.29, 180584, 66.8),
    @("NV", 382833, 1704780, 342943, 453738, 1321947, 345.31, 110795, 32.31),
    @("OH", 1510040, 2640342, 981307, 1611784, 1130302, 74.85, 630477, 64.25),
    @("OR", 528088, 486815, 528018, 486815, -41273, -7.82, -41203, -7.8),
    @("RI", 33788, 233604, 13008, 58703, 199816, 591.38, 45695, 351.28),
    @("TN", 1153458, 1151475, 1153431, 1151475, -1983, -0.17, -1956, -0.17),
    @("TX", 2115050, 5443098, 2092390, 5443098, 3328048, 157.35, 3350708, 160.14),
    @("UT", 1016921, 213235, 244846, 213235, -803686, -79.03, -31611, -12.91),
    @("VA", 423179, 2086280, 289561, 1638921, 1663101, 393, 1349360, 466),
    @("VT", 65243, 425197, 43917, 166913, 359954, 551.71, 122996, 280.06),
    @("WA", 2095234, 1161480, 785744, 1049315, -933754, -44.57, 263571, 33.54),
    @("WI", 522199, 1376055, 464968, 751288, 853856, 163.51, 286320, 61.58),
    @("WV", 66935, 179986, 60078, 139791, 113051, 168.9, 79713, 132.68),
    @("WY", 59810, 108497, 45161, 86826, 48687, 81.4, 41665, 92.26),
    @("CT", $null, 627481, $null, 411574, $null, $null, $null, $null),
    @("OK", $null, 326090, $null, 0, $null, $null, $null, $null),
    @("PA", $null, 2913035, $null, 1324715, $null, $null, $null, $null),
    @("SD", $null, 177966, $null, 148012, $null, $null, $null, $null),
)

# Clear out the old table (it spanned A1:I41) before laying down the
# new, re-ordered data so no stale cells are left behind beyond the
# new table's extent.
$ws.Range("A1:I41").ClearContents()

$r = 1
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
